$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the last-modified date in C1 ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- "MCF" sheet: bump several capacity-factor inputs to 1 ---
$mcf = $wb.Worksheets.Item("MCF")

$mcf.Range("B2").Value = 1
$mcf.Range("B3").Value = 1
$mcf.Range("B4").Value = 1
$mcf.Range("B6").Value = 1
$mcf.Range("B10").Value = 1
$mcf.Range("B11").Value = 1
$mcf.Range("B12").Value = 1
$mcf.Range("B13").Value = 1
$mcf.Range("B14").Value = 1
$mcf.Range("B16").Value = 1
$mcf.Range("B17").Value = 1
$mcf.Range("B18").Value = 1

# Move the active selection on the MCF sheet to B17, matching the saved view state.
$mcf.Activate()
$mcf.Range("B17").Select()
